$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update identifier values for rows 3-10 to match the final prototype
$ws.Range("B3").Value = "FAMILY_AFFECTED_QUESTIONS"
$ws.Range("B4").Value = "FAMILY_BREAST_QUESTIONS"
$ws.Range("B5").Value = "FAMILY_OVARIAN_QUESTIONS"
$ws.Range("B6").Value = "FAMILY_AFFECTED_GRANDMOTHER_QUESTIONS"
$ws.Range("B7").Value = "FAMILY_AFFECTED_AUNT_QUESTIONS"
$ws.Range("B8").Value = "FAMILY_AFFECTED_NIECE_QUESTIONS"
$ws.Range("B9").Value = "FAMILY_AFFECTED_HALF_SISTER_QUESTIONS"
$ws.Range("B10").Value = "PERSONAL_HISTORY_QUESTIONS"

# Remove the now-obsolete rows 11-14 entirely
$ws.Rows("11:14").Delete()

# Update selection to match final state
$ws.Range("A11").Select()
